$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 311476.12
$ws.Range("I80").Value = 381.8
$ws.Range("K80").Value = 1145.4
$ws.Range("M80").Value = -147.4000000000001
$ws.Range("H82").Value = 28580084
$ws.Range("I82").Value = 2250
$ws.Range("J82").Value = 40011220
$ws.Range("K82").Value = 6750
$ws.Range("L82").Value = 120033660
$ws.Range("M82").Value = -6344
$ws.Range("N82").Value = -120034472
$ws.Range("H83").Value = 311476.12
$ws.Range("I83").Value = 381.8
$ws.Range("K83").Value = 3436.2
$ws.Range("M83").Value = 1555.8
$ws.Range("H85").Value = 28580084
$ws.Range("I85").Value = 2250
$ws.Range("J85").Value = 40011220
$ws.Range("K85").Value = 6750
$ws.Range("L85").Value = 120033660
$ws.Range("M85").Value = -5346
$ws.Range("N85").Value = -120036468
$ws.Range("H88").Value = 883556.1
$ws.Range("I88").Value = 3625
$ws.Range("J88").Value = 1134965
$ws.Range("K88").Value = 3625
$ws.Range("L88").Value = 1134965
$ws.Range("M88").Value = -3219
$ws.Range("N88").Value = -1135777
$ws.Range("H91").Value = 883556.1
$ws.Range("I91").Value = 3625
$ws.Range("J91").Value = 1134965
$ws.Range("K91").Value = 3625
$ws.Range("L91").Value = 1134965
$ws.Range("M91").Value = -2221
$ws.Range("N91").Value = -1137773
$ws.Range("H137").Value = 1320.8
$ws.Range("I137").Value = 1086
$ws.Range("K137").Value = 3258
$ws.Range("M137").Value = -708

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1698044
$ws.Range("I32").Value = 4774.277
$ws.Range("K32").Value = 4774.277
$ws.Range("M32").Value = -4487.277
$ws.Range("H61").Value = 1468.1082
$ws.Range("I61").Value = 1386.3334
$ws.Range("J61").Value = 1619.0769
$ws.Range("K61").Value = 1386.3334
$ws.Range("L61").Value = 1619.0769
$ws.Range("M61").Value = -1174.3334
$ws.Range("N61").Value = -2043.0769
$ws.Range("H74").Value = 938.0909
$ws.Range("I74").Value = 903.4545000000001
$ws.Range("J74").Value = 1007.36365
$ws.Range("K74").Value = 903.4545000000001
$ws.Range("L74").Value = 1007.36365
$ws.Range("M74").Value = -29.45450000000005
$ws.Range("N74").Value = -2755.36365
$ws.Range("H77").Value = 938.0909
$ws.Range("I77").Value = 903.4545000000001
$ws.Range("J77").Value = 1007.36365
$ws.Range("K77").Value = 4517.2725
$ws.Range("L77").Value = 5036.81825
$ws.Range("M77").Value = -149.2725
$ws.Range("N77").Value = -13772.81825
$ws.Range("H132").Value = 2308.1875
$ws.Range("I132").Value = 1506
$ws.Range("J132").Value = 2789.5
$ws.Range("K132").Value = 4518
$ws.Range("L132").Value = 8368.5
$ws.Range("M132").Value = -1988
$ws.Range("N132").Value = -13428.5
$ws.Range("H136").Value = 1468.1082
$ws.Range("I136").Value = 1386.3334
$ws.Range("J136").Value = 1619.0769
$ws.Range("K136").Value = 4159.0002
$ws.Range("L136").Value = 4857.2307
$ws.Range("M136").Value = -1609.0002
$ws.Range("N136").Value = -9957.2307
$ws.Range("H139").Value = 43697.5
$ws.Range("J139").Value = 43697.5
$ws.Range("L139").Value = 43697.5
$ws.Range("N139").Value = -53977.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = $null
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2841.238
$ws.Range("I134").Value = 2158.9285
$ws.Range("J134").Value = 4205.857
$ws.Range("K134").Value = 6476.7855
$ws.Range("L134").Value = 12617.571
$ws.Range("M134").Value = -3941.7855
$ws.Range("N134").Value = -17687.571
$ws.Range("H140").Value = 51394.547
$ws.Range("J140").Value = 51394.547
$ws.Range("L140").Value = 51394.547
$ws.Range("N140").Value = -61754.547

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 95418.09
$ws.Range("I63").Value = 4900
$ws.Range("J63").Value = 104469.9
$ws.Range("K63").Value = 14700
$ws.Range("L63").Value = 313409.7
$ws.Range("M63").Value = -13951
$ws.Range("N63").Value = -314907.7
$ws.Range("H64").Value = 201579.8
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 251724.75
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 755174.25
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -755714.25
$ws.Range("H66").Value = 95418.09
$ws.Range("I66").Value = 4900
$ws.Range("J66").Value = 104469.9
$ws.Range("K66").Value = 44100
$ws.Range("L66").Value = 940229.1
$ws.Range("M66").Value = -40356
$ws.Range("N66").Value = -947717.1
$ws.Range("H67").Value = 201579.8
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 251724.75
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 755174.25
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -757046.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4378.8335
$ws.Range("I70").Value = 4159.4644
$ws.Range("K70").Value = 4159.4644
$ws.Range("M70").Value = -3889.4644
$ws.Range("H73").Value = 4378.8335
$ws.Range("I73").Value = 4159.4644
$ws.Range("K73").Value = 4159.4644
$ws.Range("M73").Value = -3223.4644
$ws.Range("H138").Value = 34061.8
$ws.Range("J138").Value = 34061.8
$ws.Range("L138").Value = 34061.8
$ws.Range("N138").Value = -44341.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4593.75
$ws.Range("I122").Value = 5975.3
$ws.Range("K122").Value = 17925.9
$ws.Range("M122").Value = -15475.9
$ws.Range("H132").Value = 3140.7827
$ws.Range("I132").Value = 2861.6667
$ws.Range("J132").Value = 3445.2727
$ws.Range("K132").Value = 8585.000100000001
$ws.Range("L132").Value = 10335.8181
$ws.Range("M132").Value = -6055.000100000001
$ws.Range("N132").Value = -15395.8181
$ws.Range("H136").Value = 3377.15
$ws.Range("I136").Value = 1450.8667
$ws.Range("K136").Value = 4352.6001
$ws.Range("M136").Value = -1802.6001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 710.2353000000001
$ws.Range("I107").Value = 789.125
$ws.Range("J107").Value = 520.9
$ws.Range("K107").Value = 2367.375
$ws.Range("L107").Value = 1562.7
$ws.Range("M107").Value = -447.375
$ws.Range("N107").Value = -5402.7
$ws.Range("H126").Value = 866.6667
$ws.Range("I126").Value = 900
$ws.Range("J126").Value = 800
$ws.Range("K126").Value = 2700
$ws.Range("L126").Value = 2400
$ws.Range("M126").Value = -230
$ws.Range("N126").Value = -7340
$ws.Range("H132").Value = 2854.9375
$ws.Range("I132").Value = 3166.5
$ws.Range("K132").Value = 9499.5
$ws.Range("M132").Value = -6969.5
$ws.Range("H136").Value = 1759.0769
$ws.Range("I136").Value = 1739
$ws.Range("K136").Value = 5217
$ws.Range("M136").Value = -2667
